$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header: "note" label moves from D1's old string slot, text unchanged
$ws.Range("D1").Value = "note"

# Row 4: supportPlateThickness -> height, value 0.25 -> 150, unit in -> mm
# (edited first, so its label lands earliest among the new shared strings)
$ws.Range("A4").Value = "height"
$ws.Range("B4").NumberFormat = "0.000"
$ws.Range("B4").Value = 150
$ws.Range("C4").Value = "mm"

# Row 2: plateWidth -> width, unit in -> mm (value unchanged)
$ws.Range("A2").Value = "width"
$ws.Range("B2").Value = 300
$ws.Range("C2").Value = "mm"

# Row 3: plateDepth -> depth, unit in -> mm (value unchanged)
$ws.Range("A3").Value = "depth"
$ws.Range("B3").Value = 1000
$ws.Range("C3").Value = "mm"

# Row 5: supportPlateWidth -> suppourtPlateLength, value 50 -> 100, format 0.00000 -> 0.000, unit mm
$ws.Range("A5").Value = "suppourtPlateLength"
$ws.Range("B5").NumberFormat = "0.000"
$ws.Range("B5").Value = 100
$ws.Range("C5").Value = "mm"

# Row 6: supportPlateFillet -> acrylicThickness, value 10 -> 3.125, unit mm
$ws.Range("A6").Value = "acrylicThickness"
$ws.Range("B6").Value = 3.125
$ws.Range("C6").Value = "mm"

# Rows 7-9: clear out the old separation / numPanels / centerVerticalRatio rows entirely
$ws.Range("A7").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()

$ws.Range("A8").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()

$ws.Range("A9").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("C9").ClearContents()

# Selection moves to B7
$ws.Range("B7").Select()
